# iron_native / Word COM-interop script
# Applies the template changes described by the commit:
#   "feat: buat permohonan - menambahkan form untuk masing-masing surat"
#
# Summary of changes to src/assets/docx/templete/4.docx:
#   1. Remove the stray _GoBack bookmark that sat after the {agama} run.
#   2. In the signature table cell: keep the two existing blank lines,
#      add a new bold/centered blank line, then a centered line that
#      holds the {nama} placeholder, and blank out the old bold dotted
#      signature-line text (its paragraph mark / formatting stays).
#   3. Re-add the _GoBack bookmark on the (now) first blank paragraph
#      following the table.
#
# We use Range.InsertXML with hand-built OOXML fragments (rather than
# the generic paragraph-insertion methods) so every <w:pPr>/<w:rPr>
# comes out byte-for-byte as in the target, instead of relying on
# whatever default formatting a blank InsertParagraphBefore() would
# invent.

$d = $word.ActiveDocument

$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------
# 1) Second blank paragraph right after the table (currently #28):
#    add back the _GoBack bookmark pair inside its (otherwise
#    untouched) pPr.
# ---------------------------------------------------------------------
$pBookmark = $d.Paragraphs.Item(28)
$bookmarkXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3686"/><w:tab w:val="left" w:pos="3969"/></w:tabs><w:ind w:firstLine="709"/><w:jc w:val="center"/><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$pBookmark.Range.InsertXML($bookmarkXml)

# ---------------------------------------------------------------------
# 2) Dotted signature-line paragraph in the table (#26): drop the run
#    with the dots, leave the (bold, centered) paragraph mark in place.
# ---------------------------------------------------------------------
$pDots = $d.Paragraphs.Item(26)
$dotsXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3686"/><w:tab w:val="left" w:pos="3969"/></w:tabs><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p>
'@
$pDots.Range.InsertXML($dotsXml)

# ---------------------------------------------------------------------
# 3) Second (of two identical) blank paragraphs in the table (#25):
#    split it into a new bold/centered blank paragraph followed by a
#    centered paragraph carrying the {nama} placeholder.
# ---------------------------------------------------------------------
$pSplit = $d.Paragraphs.Item(25)
$splitXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3686"/><w:tab w:val="left" w:pos="3969"/></w:tabs><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:tabs><w:tab w:val="left" w:pos="3686"/><w:tab w:val="left" w:pos="3969"/></w:tabs><w:spacing w:after="0" w:line="360" w:lineRule="auto"/><w:jc w:val="center"/><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>{nama}</w:t></w:r></w:p>
'@
$pSplit.Range.InsertXML($splitXml)

# ---------------------------------------------------------------------
# 4) "Agama" paragraph (#9): drop the leftover _GoBack bookmark that
#    used to sit right after the {agama} run.
# ---------------------------------------------------------------------
$pAgama = $d.Paragraphs.Item(9)
$agamaXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="5"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="left" w:pos="3686"/><w:tab w:val="left" w:pos="3969"/></w:tabs><w:spacing w:line="480" w:lineRule="auto"/><w:ind w:left="851" w:hanging="425"/><w:jc w:val="both"/><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Agama</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>:</w:t></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:hint="default" w:ascii="Cambria" w:hAnsi="Cambria" w:cs="Cambria"/><w:sz w:val="24"/><w:szCs w:val="24"/><w:lang w:val="en-GB"/></w:rPr><w:t>{agama}</w:t></w:r></w:p>
'@
$pAgama.Range.InsertXML($agamaXml)

Write-Output ("Paragraphs after edit: {0}" -f $d.Paragraphs.Count)
